# ---------------------------------------------------------------------------
# mxif-metadata.xlsx: add a "version" / "description" pair of leading columns
# to the "Export as TSV" sheet, plus a new "version list" lookup sheet that
# backs a data-validation list on the new "version" column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 0) Snapshot the existing header-row cell comments (A1:Z1) before we start
#    moving columns around, so we can re-create them in their new homes.
# ---------------------------------------------------------------------------
$commentTexts = @()
for ($c = 1; $c -le 26; $c++) {
    $cell = $ws1.Cells.Item(1, $c)
    $cm = $cell.Comment
    $commentTexts += $cm.Text()
}

# Delete the old comments now (their cells are about to shift two columns
# right, but comments do not follow automatically).
for ($c = 1; $c -le 26; $c++) {
    $cell = $ws1.Cells.Item(1, $c)
    $cell.Comment.Delete()
}

# ---------------------------------------------------------------------------
# 1) Add the new "version list" sheet right after "Export as TSV".
# ---------------------------------------------------------------------------
$versionSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").Value = "1"

# ---------------------------------------------------------------------------
# 2) Insert two new leading columns (A, B) on the main sheet. Everything
#    else (headers, data validations) shifts two columns to the right.
# ---------------------------------------------------------------------------
$ws1.Range("A1:B1").EntireColumn.Insert()

$ws1.Cells.Item(1, 1).Value = "version"
$ws1.Cells.Item(1, 2).Value = "description"

# ---------------------------------------------------------------------------
# 3) Re-create the header comments, shifted two columns to the right, plus
#    the two new ones for "version" and "description".
# ---------------------------------------------------------------------------
$ws1.Cells.Item(1, 1).AddComment("Version of the schema to use when validating this metadata.")
$ws1.Cells.Item(1, 2).AddComment("Free-text description of this assay.")

for ($c = 1; $c -le 26; $c++) {
    $cell = $ws1.Cells.Item(1, $c + 2)
    $cell.AddComment($commentTexts[$c - 1])
}

# ---------------------------------------------------------------------------
# 4) Add the data validation list for the new "version" column.
# ---------------------------------------------------------------------------
$verRange = $ws1.Range("A2:A1048576")
$verRange.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$verRange.Validation.ErrorTitle = "Value must come from list"
$verRange.Validation.ErrorMessage = "Value must be one of: 1."
$verRange.Validation.ShowInput = $true
$verRange.Validation.ShowError = $true
